$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same values applied to every data row, B2:Q26)
# NOTE: scientific-notation literals (1e-05) are not accepted by the
# script parser, so very small/large values are written as plain
# fixed-point decimals that round-trip to the exact same IEEE-754 double.
$values = @{
    "B" = 0.9999845286516351
    "C" = 0.9991389381953254
    "D" = 0.9999999763723475
    "E" = 0.9999913153482848
    "F" = 0.9999958442422094
    "G" = 0.00001444182458493536
    "H" = 0.000803763398418822
    "I" = 0.000000007103075553283338
    "J" = 0.00000391546207290299
    "K" = 0.000001961282555009019
    "L" = 0.0001971186733916371
    "M" = 0.003800240069381849
    "N" = 1.00001125188972
    "O" = 0.003962024156716867
    "P" = 136.2907641523688
    "Q" = 205.7666861698563
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $range = $ws.Range($col + "2:" + $col + "26")
    $range.Value = $val
}
